$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.369.58"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "2.604.34"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.43"
$ws.Range("E5").Value = "  +3.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.15"
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "2.612.61"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.156"
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.371"
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").Value = "3.058.61"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.42"
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").Value = "60.370.59"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("D18").Value = "2.608.92"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.56"
$ws.Range("E19").Value = "  +10.18%  "
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "347.91"
$ws.Range("E21").Value = "  +2.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.90"
$ws.Range("E22").Value = "  +3.87%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.524"
$ws.Range("E24").Value = "  +5.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.11"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.09"
$ws.Range("E28").Value = "  +7.94%  "
$ws.Range("D29").Value = "0.0₃0796"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("E30").Value = "  +9.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.41"
$ws.Range("E31").Value = "  +4.43%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "164.72"
$ws.Range("E33").Value = "  +3.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.46"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.33"
$ws.Range("E35").Value = "  +6.32%  "
$ws.Range("E36").Value = "  +9.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.989"
$ws.Range("E37").Value = "  +8.60%  "
$ws.Range("E38").Value = "  +8.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.13"
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "314.80"
$ws.Range("E40").Value = "  +7.79%  "
$ws.Range("E41").Value = "  +6.09%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "134.85"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0994"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.07"
$ws.Range("E46").Value = "  +6.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.88"
$ws.Range("E47").Value = "  +4.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.607"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0552"
$ws.Range("E49").Value = "  +2.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.09"
$ws.Range("E50").Value = "  +6.85%  "
$ws.Range("E51").Value = "  +1.37%  "
